$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "TestScenario_1.TestCase_1"
$ws.Range("C2").Value = "New Opportunity"
$ws.Range("D2").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "Step 1"
$ws.Range("G2").Value = "Click on the Opportunity tab, and click on New button"
$ws.Range("H2").Value = "User should be navigated to the New  Opportunity Page"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# Row 3
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Valid value for required field Amount "
$ws.Range("F3").Value = "Step 2"
$ws.Range("G3").Value = "Input valid value in the  Amount field."
$ws.Range("H3").Value = "User should be able to input value for the Amount field."
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""

# Row 4
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Valid value for required field Close Date "
$ws.Range("F4").Value = "Step 3"
$ws.Range("G4").Value = "Input valid value in the  Close Date field."
$ws.Range("H4").Value = "User should be able to input value for the Close Date field."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""

# Row 5
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Valid value for required field Opportunity Name "
$ws.Range("F5").Value = "Step 4"
$ws.Range("G5").Value = "Input valid value in the  Opportunity Name field."
$ws.Range("H5").Value = "User should be able to input value for the Opportunity Name field."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""

# Row 6
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "Valid value for required field  "
$ws.Range("F6").Value = "Step 5"
$ws.Range("G6").Value = "Input valid value in the   field."
$ws.Range("H6").Value = "User should be able to input value for the  field."
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""

# Row 7
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "Valid value for required field Stage, value should be equals Proposal/Price Quote to submit for Approval towards to the assigned approver "
$ws.Range("F7").Value = "Step 6"
$ws.Range("G7").Value = "Input valid value in the  Stage field."
$ws.Range("H7").Value = "Value accepted for Stage field."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""

# Row 8
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "Step 7"
$ws.Range("G8").Value = "Click on Save button to save Opportunity with fields"
$ws.Range("H8").Value = "User should be able to validate that a New Opportunity is created"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""

# Row 9
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "Step 8"
$ws.Range("G9").Value = "On the Opportunity page, Click on 'Submit for Approval' button to Submit for Approval."
$ws.Range("H9").Value = "Pop-up confirming to submit the record for Approval is displayed."
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""

# Row 10
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "Step 9"
$ws.Range("G10").Value = "Click on 'Cancel' button to prevent submission for approval."
$ws.Range("H10").Value = "User is redirected to the Opportunity Page"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""

# Row 11
$ws.Range("A11").Value = "TestScenario_1"
$ws.Range("B11").Value = "TestScenario_1.TestCase_1"
$ws.Range("C11").Value = "Approve/Reject Opportunity"
$ws.Range("D11").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "Step 1"
$ws.Range("G11").Value = "Click on the Opportunity tab"
$ws.Range("H11").Value = "User should be navigated to the Opportunity Page"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""

# Row 12
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "Step 2"
$ws.Range("G12").Value = "From the list of the  Opportunitys displayed, select the appropriate Opportunity."
$ws.Range("H12").Value = "User should be navigated to the Opportunity details page."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""

# Row 13
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "Step 3"
$ws.Range("G13").Value = "Scroll down the Opportunity page to locate the 'Approval History' section."
$ws.Range("H13").Value = "User should be able to view the pending Approve/Reject requests listed."
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""

# Row 14
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = "Step 4"
$ws.Range("G14").Value = "To Approve/Reject the record's request, Click on 'Approve/Reject' link."
$ws.Range("H14").Value = "User should be navigated to the Approval Request, Account:(Approver Name) Page."
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""

# Row 15
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = "Step 5"
$ws.Range("G15").Value = "In the 'Approve/Reject Approval Request' section, input comments if required."
$ws.Range("H15").Value = "User should be able to input appropriate comments."
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""

# Row 16
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = "Step 6"
$ws.Range("G16").Value = "Click on the 'Approve' or the 'Reject' button to either Approve or Reject the request."
$ws.Range("H16").Value = "User should be able to either 'Approve' or 'Reject' the request."
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""

# Row 17
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = "Step 7"
$ws.Range("G17").Value = "On performing either 'Approval' or 'Rejection' action, user is navigated to the Opportunity request page."
$ws.Range("H17").Value = "User should be navigated to the Opportunity details page."
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""

# Row 18
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = "Step 8"
$ws.Range("G18").Value = "Scroll down the Opportunity page to locate the 'Approval History' section."
$ws.Range("H18").Value = "User should be able to view the request is either 'Approved' or 'Rejected'."
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""

# Expand column E width to match new content (best achievable via ColumnWidth COM property)
$ws.Columns.Item(5).ColumnWidth = 124

# Resize table / ListObject and autofilter range to cover the new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J18"))
